$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2464.9312
$ws.Range("J112").Value = 2912.348
$ws.Range("L112").Value = 8737.044
$ws.Range("N112").Value = -10953.044
$ws.Range("H116").Value = 2360.8708
$ws.Range("I116").Value = 2257.524
$ws.Range("J116").Value = 2577.9
$ws.Range("K116").Value = 2257.524
$ws.Range("L116").Value = 2577.9
$ws.Range("M116").Value = 1184.476
$ws.Range("N116").Value = -9461.9
$ws.Range("H129").Value = 822.0513
$ws.Range("J129").Value = 1054.7407
$ws.Range("L129").Value = 3164.2221
$ws.Range("N129").Value = -13164.2221
$ws.Range("H132").Value = 7098915
$ws.Range("I132").Value = 9808600
$ws.Range("J132").Value = 12047.538
$ws.Range("K132").Value = 29425800
$ws.Range("L132").Value = 36142.614
$ws.Range("M132").Value = -29423270
$ws.Range("N132").Value = -41202.614
$ws.Range("H137").Value = 1158.6383
$ws.Range("I137").Value = 832.9048
$ws.Range("K137").Value = 2498.7144
$ws.Range("M137").Value = 51.28560000000016
$ws.Range("H138").Value = 1522.77
$ws.Range("I138").Value = 997.0714
$ws.Range("J138").Value = 1608.3489
$ws.Range("K138").Value = 2991.2142
$ws.Range("L138").Value = 4825.0467
$ws.Range("M138").Value = 2148.7858
$ws.Range("N138").Value = -15105.0467
$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3721.123
$ws.Range("I32").Value = 3297.9
$ws.Range("K32").Value = 3297.9
$ws.Range("M32").Value = -3010.9
$ws.Range("H132").Value = 1158.6508
$ws.Range("I132").Value = 885.6731
$ws.Range("K132").Value = 2657.0193
$ws.Range("M132").Value = -127.0192999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4205.8
$ws.Range("I86").Value = 4205.8
$ws.Range("K86").Value = 4205.8
$ws.Range("M86").Value = -3082.8
$ws.Range("H89").Value = 4205.8
$ws.Range("I89").Value = 4205.8
$ws.Range("K89").Value = 21029
$ws.Range("M89").Value = -15413
$ws.Range("H134").Value = 4160.9
$ws.Range("I134").Value = 943.51514
$ws.Range("J134").Value = 19328.572
$ws.Range("K134").Value = 2830.54542
$ws.Range("L134").Value = 57985.716
$ws.Range("M134").Value = -295.5454199999999
$ws.Range("N134").Value = -63055.716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2082.4
$ws.Range("I31").Value = 2073.8
$ws.Range("K31").Value = 2073.8
$ws.Range("M31").Value = -1778.8
$ws.Range("H34").Value = 2082.4
$ws.Range("I34").Value = 2073.8
$ws.Range("K34").Value = 2073.8
$ws.Range("M34").Value = -1871.8
$ws.Range("H58").Value = 613.40424
$ws.Range("I58").Value = 671.2917
$ws.Range("J58").Value = 553
$ws.Range("K58").Value = 671.2917
$ws.Range("L58").Value = 553
$ws.Range("M58").Value = -468.2917
$ws.Range("N58").Value = -959
$ws.Range("H134").Value = 612.9375
$ws.Range("I134").Value = 586.3611
$ws.Range("J134").Value = 692.6667
$ws.Range("K134").Value = 1759.0833
$ws.Range("L134").Value = 2078.0001
$ws.Range("M134").Value = 775.9167000000002
$ws.Range("N134").Value = -7148.0001
$ws.Range("H136").Value = 613.40424
$ws.Range("I136").Value = 671.2917
$ws.Range("J136").Value = 553
$ws.Range("K136").Value = 2013.8751
$ws.Range("L136").Value = 1659
$ws.Range("M136").Value = 536.1249
$ws.Range("N136").Value = -6759

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1192.4865
$ws.Range("I5").Value = 1351.4073
$ws.Range("J5").Value = 763.4
$ws.Range("K5").Value = 4054.2219
$ws.Range("L5").Value = 2290.2
$ws.Range("M5").Value = -3942.2219
$ws.Range("N5").Value = -2514.2
$ws.Range("H135").Value = 1192.4865
$ws.Range("I135").Value = 1351.4073
$ws.Range("J135").Value = 763.4
$ws.Range("K135").Value = 12162.6657
$ws.Range("L135").Value = 6870.599999999999
$ws.Range("M135").Value = -9627.665700000001
$ws.Range("N135").Value = -11940.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 12000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 12000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 12000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -12278
$ws.Range("H17").Value = 2750
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2750
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2750
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3086
$ws.Range("H23").Value = 1990
$ws.Range("I23").Value = 1990
$ws.Range("K23").Value = 1990
$ws.Range("M23").Value = -1767
$ws.Range("H34").Value = 23000
$ws.Range("J34").Value = 23000
$ws.Range("L34").Value = 23000
$ws.Range("N34").Value = -23536
$ws.Range("H74").Value = 51533.332
$ws.Range("J74").Value = 51533.332
$ws.Range("L74").Value = 51533.332
$ws.Range("N74").Value = -53405.332
$ws.Range("H76").Value = 23000
$ws.Range("J76").Value = 23000
$ws.Range("L76").Value = 23000
$ws.Range("N76").Value = -23630
$ws.Range("H77").Value = 51533.332
$ws.Range("J77").Value = 51533.332
$ws.Range("L77").Value = 154599.996
$ws.Range("N77").Value = -163959.996
$ws.Range("H79").Value = 23000
$ws.Range("J79").Value = 23000
$ws.Range("L79").Value = 23000
$ws.Range("N79").Value = -25184
$ws.Range("H113").Value = 1451.5
$ws.Range("I113").Value = 1414.9333
$ws.Range("K113").Value = 1414.9333
$ws.Range("M113").Value = 755.0667000000001
$ws.Range("H132").Value = 1794.7059
$ws.Range("I132").Value = 1321.3572
$ws.Range("K132").Value = 3964.0716
$ws.Range("M132").Value = -1434.0716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1667.6666
$ws.Range("J22").Value = 2001
$ws.Range("L22").Value = 2001
$ws.Range("N22").Value = -2591
$ws.Range("H27").Value = 1667.6666
$ws.Range("J27").Value = 2001
$ws.Range("L27").Value = 2001
$ws.Range("N27").Value = -2215
$ws.Range("H61").Value = 2046
$ws.Range("I61").Value = 1398.1666
$ws.Range("K61").Value = 1398.1666
$ws.Range("M61").Value = -1196.1666
$ws.Range("H82").Value = 1657.9615
$ws.Range("I82").Value = 1645.6364
$ws.Range("K82").Value = 1645.6364
$ws.Range("M82").Value = -1284.6364
$ws.Range("H85").Value = 1657.9615
$ws.Range("I85").Value = 1645.6364
$ws.Range("K85").Value = 1645.6364
$ws.Range("M85").Value = -397.6364000000001
$ws.Range("H113").Value = 2046
$ws.Range("I113").Value = 1398.1666
$ws.Range("K113").Value = 1398.1666
$ws.Range("M113").Value = 771.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 525.86365
$ws.Range("I136").Value = 258.14285
$ws.Range("K136").Value = 774.4285500000001
$ws.Range("M136").Value = 1775.57145
